$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 93
$ws.Range("H93").Value = 40733.668
$ws.Range("J93").Value = 40733.668
$ws.Range("L93").Value = 40733.668
$ws.Range("N93").Value = -45725.668
# Row 124
$ws.Range("H124").Value = 41772.727
$ws.Range("J124").Value = 41772.727
$ws.Range("L124").Value = 41772.727
$ws.Range("N124").Value = -51592.727
# Row 126
$ws.Range("H126").Value = 47800
$ws.Range("J126").Value = 47800
$ws.Range("L126").Value = 47800
$ws.Range("N126").Value = -57680
# Row 128
$ws.Range("H128").Value = 37566.668
$ws.Range("J128").Value = 37566.668
$ws.Range("L128").Value = 37566.668
$ws.Range("N128").Value = -47526.668
# Row 130
$ws.Range("H130").Value = 48532.082
$ws.Range("J130").Value = 48532.082
$ws.Range("L130").Value = 48532.082
$ws.Range("N130").Value = -58572.082
# Row 133
$ws.Range("H133").Value = 51246.668
$ws.Range("J133").Value = 51246.668
$ws.Range("L133").Value = 51246.668
$ws.Range("N133").Value = -61366.668

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 44
$ws.Range("H44").Value = 20212
$ws.Range("J44").Value = 20212
$ws.Range("L44").Value = 20212
$ws.Range("N44").Value = -21188
# Row 80
$ws.Range("H80").Value = 33688.57
$ws.Range("J80").Value = 33688.57
$ws.Range("L80").Value = 33688.57
$ws.Range("N80").Value = -35684.57
# Row 83
$ws.Range("H83").Value = 33688.57
$ws.Range("J83").Value = 33688.57
$ws.Range("L83").Value = 101065.71
$ws.Range("N83").Value = -111049.71
# Row 109
$ws.Range("H109").Value = 33699.332
$ws.Range("J109").Value = 33699.332
$ws.Range("L109").Value = 33699.332
$ws.Range("N109").Value = -36473.332
# Row 122
$ws.Range("H122").Value = 1665.2084
$ws.Range("I122").Value = 1499.1765
$ws.Range("J122").Value = 2068.4285
$ws.Range("K122").Value = 4497.529500000001
$ws.Range("L122").Value = 6205.2855
$ws.Range("M122").Value = -2047.529500000001
$ws.Range("N122").Value = -11105.2855
# Row 125
$ws.Range("H125").Value = 112548600
$ws.Range("J125").Value = 112548600
$ws.Range("L125").Value = 112548600
$ws.Range("N125").Value = -112558440
# Row 127
$ws.Range("H127").Value = 54572
$ws.Range("J127").Value = 54572
$ws.Range("L127").Value = 54572
$ws.Range("N127").Value = -64492
# Row 129
$ws.Range("H129").Value = 49249.668
$ws.Range("J129").Value = 49249.668
$ws.Range("L129").Value = 49249.668
$ws.Range("N129").Value = -59249.668
# Row 130
$ws.Range("H130").Value = 27143
$ws.Range("J130").Value = 27143
$ws.Range("L130").Value = 27143
$ws.Range("N130").Value = -37183
# Row 131
$ws.Range("H131").Value = 59882.668
$ws.Range("J131").Value = 59882.668
$ws.Range("L131").Value = 59882.668
$ws.Range("N131").Value = -69962.66800000001
# Row 135
$ws.Range("H135").Value = 23306.785
$ws.Range("J135").Value = 23306.785
$ws.Range("L135").Value = 23306.785
$ws.Range("N135").Value = -33446.785

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 10430.2
$ws.Range("I26").Value = 6787.75
$ws.Range("J26").Value = 25000
$ws.Range("K26").Value = 6787.75
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = -6495.75
$ws.Range("N26").Value = -25584
# Row 35
$ws.Range("H35").Value = 23137.6
$ws.Range("I35").Value = 15000
$ws.Range("K35").Value = 15000
$ws.Range("M35").Value = -14690
# Row 82
$ws.Range("H82").Value = 16812.117
$ws.Range("J82").Value = 27649.889
$ws.Range("L82").Value = 27649.889
$ws.Range("N82").Value = -28415.889
# Row 85
$ws.Range("H85").Value = 16812.117
$ws.Range("J85").Value = 27649.889
$ws.Range("L85").Value = 27649.889
$ws.Range("N85").Value = -30301.889
# Row 122
$ws.Range("H122").Value = 54124
$ws.Range("J122").Value = 54124
$ws.Range("L122").Value = 54124
$ws.Range("N122").Value = -63924
# Row 124
$ws.Range("H124").Value = 38755.555
$ws.Range("J124").Value = 38755.555
$ws.Range("L124").Value = 38755.555
$ws.Range("N124").Value = -48575.555
# Row 125
$ws.Range("H125").Value = 51788.89
$ws.Range("J125").Value = 51788.89
$ws.Range("L125").Value = 51788.89
$ws.Range("N125").Value = -61628.89
# Row 126
$ws.Range("H126").Value = 58792
$ws.Range("J126").Value = 58792
$ws.Range("L126").Value = 58792
$ws.Range("N126").Value = -68672
# Row 129
$ws.Range("H129").Value = 49992.668
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 49992.668
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 49992.668
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -59992.668
# Row 130
$ws.Range("H130").Value = 59968
$ws.Range("J130").Value = 59968
$ws.Range("L130").Value = 59968
$ws.Range("N130").Value = -70008
# Row 135
$ws.Range("H135").Value = 40248.57
$ws.Range("I135").Value = 10000
$ws.Range("J135").Value = 62935
$ws.Range("K135").Value = 10000
$ws.Range("L135").Value = 62935
$ws.Range("N135").Value = -73075
$ws.Range("M135").Value = -4930

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 20
$ws.Range("H20").Value = 52489.5
$ws.Range("J20").Value = 52489.5
$ws.Range("L20").Value = 52489.5
$ws.Range("N20").Value = -52961.5
# Row 30
$ws.Range("H30").Value = 52489.5
$ws.Range("J30").Value = 52489.5
$ws.Range("L30").Value = 52489.5
$ws.Range("N30").Value = -52671.5
# Row 75
$ws.Range("H75").Value = 43800
$ws.Range("J75").Value = 43800
$ws.Range("L75").Value = 43800
$ws.Range("N75").Value = -45796
# Row 78
$ws.Range("H78").Value = 43800
$ws.Range("J78").Value = 43800
$ws.Range("L78").Value = 131400
$ws.Range("N78").Value = -141384
# Row 97
$ws.Range("H97").Value = 9890
$ws.Range("J97").Value = 9890
$ws.Range("L97").Value = 9890
$ws.Range("N97").Value = -11872
# Row 127
$ws.Range("H127").Value = 60000
$ws.Range("J127").Value = 60000
$ws.Range("L127").Value = 60000
$ws.Range("N127").Value = -69920
# Row 128
$ws.Range("H128").Value = 52489.5
$ws.Range("J128").Value = 52489.5
$ws.Range("L128").Value = 52489.5
$ws.Range("N128").Value = -62449.5
# Row 135
$ws.Range("H135").Value = 42728.57
$ws.Range("J135").Value = 42728.57
$ws.Range("L135").Value = 42728.57
$ws.Range("N135").Value = -52868.57

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 1005
$ws.Range("I43").Value = 1005
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1005
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("M43").Value = -854
# Row 46
$ws.Range("H46").Value = 18767.3
$ws.Range("J46").Value = 18767.3
$ws.Range("L46").Value = 18767.3
$ws.Range("N46").Value = -19079.3
# Row 93
$ws.Range("H93").Value = 9890
$ws.Range("J93").Value = 9890
$ws.Range("L93").Value = 9890
$ws.Range("N93").Value = -13634
# Row 128
$ws.Range("H128").Value = 38575
$ws.Range("J128").Value = 38575
$ws.Range("L128").Value = 38575
$ws.Range("N128").Value = -48535
# Row 133
$ws.Range("H133").Value = 34831.11
$ws.Range("J133").Value = 34831.11
$ws.Range("L133").Value = 34831.11
$ws.Range("N133").Value = -44951.11

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 92
$ws.Range("H92").Value = 35011
$ws.Range("J92").Value = 35011
$ws.Range("L92").Value = 35011
$ws.Range("N92").Value = -40003
# Row 96
$ws.Range("H96").Value = 20197
$ws.Range("J96").Value = 20197
$ws.Range("L96").Value = 20197
$ws.Range("N96").Value = -25689
# Row 123
$ws.Range("H123").Value = 53500
$ws.Range("J123").Value = 53500
$ws.Range("L123").Value = 53500
$ws.Range("N123").Value = -63300
# Row 128
$ws.Range("H128").Value = 42756.125
$ws.Range("J128").Value = 42756.125
$ws.Range("L128").Value = 42756.125
$ws.Range("N128").Value = -52716.125
# Row 129
$ws.Range("H129").Value = 34369.6
$ws.Range("J129").Value = 34369.6
$ws.Range("L129").Value = 34369.6
$ws.Range("N129").Value = -44369.6
# Row 130
$ws.Range("H130").Value = 47478.75
$ws.Range("J130").Value = 47478.75
$ws.Range("L130").Value = 47478.75
$ws.Range("N130").Value = -57518.75
# Row 133
$ws.Range("H133").Value = 22068.5
$ws.Range("I133").Value = 15296
$ws.Range("J133").Value = 24326
$ws.Range("K133").Value = 15296
$ws.Range("L133").Value = 24326
$ws.Range("N133").Value = -29386
$ws.Range("M133").Value = -12766

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 40057
$ws.Range("J64").Value = 40057
$ws.Range("L64").Value = 40057
$ws.Range("N64").Value = -40553
# Row 67
$ws.Range("H67").Value = 40057
$ws.Range("J67").Value = 40057
$ws.Range("L67").Value = 40057
$ws.Range("N67").Value = -41773
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
# Row 108
$ws.Range("H108").Value = 35626
$ws.Range("J108").Value = 35626
$ws.Range("L108").Value = 35626
$ws.Range("N108").Value = -43306
# Row 123
$ws.Range("H123").Value = 54692
$ws.Range("J123").Value = 54692
$ws.Range("L123").Value = 54692
$ws.Range("N123").Value = -64492
# Row 128
$ws.Range("H128").Value = 45327.31
$ws.Range("J128").Value = 45327.31
$ws.Range("L128").Value = 45327.31
$ws.Range("N128").Value = -55287.31
# Row 135
$ws.Range("H135").Value = 46048.934
$ws.Range("J135").Value = 46048.934
$ws.Range("L135").Value = 46048.934
$ws.Range("N135").Value = -56188.934
